# Applies the "Apply missing Part A doc headings" change:
#  1) Adds <w:outlineLvl w:val="2"/> to the "8. Are there any arrest issues..." paragraph
#  2) Adds <w:outlineLvl w:val="2"/> to the "15. Current Risk of Serious Harm..." paragraph
#  3) Moves the "21. What alternatives to recall..." heading out of the table and into its
#     own bordered paragraph (preceded by a page-break paragraph and a spacer paragraph),
#     removing the now-empty heading row from the table.
#  4) Adds <w:outlineLvl w:val="2"/> to the "22. Select the proposed recall type..." paragraph

$d = $word.ActiveDocument

function Set-OutlineLevelByText($doc, $searchText, $level) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Text not found: " + $searchText)
    }
    $rng.ParagraphFormat.OutlineLevel = $level
}

function Find-ParagraphIndexByText($doc, $searchText) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like ("*" + $searchText + "*")) {
            return $i
        }
    }
    return -1
}

# --- 1) & 2): simple outlineLvl=2 (wdOutlineLevel 3) additions -------------
Set-OutlineLevelByText $d "Are there any arrest issues of which police should be aware?" 3
Set-OutlineLevelByText $d "Current Risk of Serious Harm Assessment at time of this recall" 3

# --- 3) Insert the new "page break + spacer + 21. heading" paragraphs ------
# Do this BEFORE touching $d.Tables so paragraph indices stay in document order
# (accessing a specific Tables item shifts how table paragraphs are indexed).
$anchorIdx = Find-ParagraphIndexByText $d "response_to_probation"
if ($anchorIdx -eq -1) {
    throw "Could not find anchor paragraph {{response_to_probation}}"
}
$anchor = $d.Paragraphs.Item($anchorIdx)
$null = $anchor.Range.InsertParagraphAfter()
$null = $anchor.Range.InsertParagraphAfter()
$null = $anchor.Range.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Item($anchorIdx + 1)
$newPara2 = $d.Paragraphs.Item($anchorIdx + 2)
$newPara3 = $d.Paragraphs.Item($anchorIdx + 3)

$xmlHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$para1Xml = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="4"/><w:szCs w:val="4"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="4"/><w:szCs w:val="4"/><w:u w:val="single"/></w:rPr><w:br w:type="page"/></w:r></w:p>'
$para2Xml = '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="4"/><w:szCs w:val="4"/><w:u w:val="single"/></w:rPr></w:pPr></w:p>'
$para3Xml = '<w:p><w:pPr><w:pBdr><w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="9" w:color="auto"/></w:pBdr><w:tabs><w:tab w:val="right" w:pos="9184"/></w:tabs><w:outlineLvl w:val="2"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="800080"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="800080"/></w:rPr><w:t xml:space="preserve">21. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="800080"/></w:rPr><w:t>What alternatives to recall have been taken to try to secure compliance and manage risk prior to requesting recall?  Provide full details below including dates:</w:t></w:r></w:p>'

$null = $newPara1.Range.InsertXML($xmlHeader + $para1Xml + $xmlFooter)
$null = $newPara2.Range.InsertXML($xmlHeader + $para2Xml + $xmlFooter)
$null = $newPara3.Range.InsertXML($xmlHeader + $para3Xml + $xmlFooter)

# --- Remove the now-redundant "21." heading row from its table -------------
$removed = $false
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Rows.Item(1).Range.Text -like "*What alternatives to recall have been taken*") {
        $null = $t.Rows.Item(1).Delete()
        $removed = $true
        break
    }
}
if (-not $removed) {
    throw "Could not find the '21.' heading row to remove"
}

# --- 4) simple outlineLvl=2 (wdOutlineLevel 3) addition ---------------------
Set-OutlineLevelByText $d "Select the proposed recall type, having considered the information above" 3

Write-Output "Applied Part A heading changes"
